$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two rows (old rows 142 and 143); remaining rows shift up,
# row count goes from 143 to 141 data+header rows (dimension A1:AC141).
$ws.Range("A142:A143").EntireRow.Delete()

# Row 138: fill in the match result (FTHG/FTAG/FTR) and update PL* columns
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 2
$ws.Range("J138").Value = "A"
$ws.Range("W138").Value = -1
$ws.Range("X138").Value = -1
$ws.Range("Y138").Value = 4.5
$ws.Range("Z138").Value = -1
$ws.Range("AA138").Value = 0.825
$ws.Range("AB138").Value = -1
$ws.Range("AC138").Value = 1.025

# Row 139: updated match data
$ws.Range("B139").Value = 6818335
$ws.Range("E139").Value = 45361.375
$ws.Range("F139").Value = "Paksi"
$ws.Range("G139").Value = "Debreceni VSC"
$ws.Range("K139").Value = 2.1
$ws.Range("L139").Value = 3.4
$ws.Range("M139").Value = 3.4
$ws.Range("N139").Value = 2
$ws.Range("O139").Value = 3.4
$ws.Range("P139").Value = 3.6
$ws.Range("Q139").Value = -0.5
$ws.Range("R139").Value = 2
$ws.Range("S139").Value = 1.85
$ws.Range("U139").Value = 1.825
$ws.Range("V139").Value = 2.025

# Row 140: updated match data
$ws.Range("B140").Value = 6818338
$ws.Range("E140").Value = 45361.47916666666
$ws.Range("F140").Value = "MOL Fehervar FC"
$ws.Range("G140").Value = "Ferencvarosi TC"
$ws.Range("K140").Value = 5.25
$ws.Range("L140").Value = 4
$ws.Range("M140").Value = 1.571
$ws.Range("N140").Value = 6
$ws.Range("O140").Value = 4.2
$ws.Range("P140").Value = 1.5
$ws.Range("Q140").Value = 1
$ws.Range("R140").Value = 2.025
$ws.Range("S140").Value = 1.825
$ws.Range("T140").Value = 2.75
$ws.Range("U140").Value = 1.875
$ws.Range("V140").Value = 1.975

# Row 141: updated match data
$ws.Range("B141").Value = 6818334
$ws.Range("E141").Value = 45361.66666666666
$ws.Range("F141").Value = "Zalaegerszegi TE"
$ws.Range("G141").Value = "MTK Budapest"
$ws.Range("K141").Value = 2.2
$ws.Range("M141").Value = 3.1
$ws.Range("N141").Value = 2.2
$ws.Range("P141").Value = 3.1
$ws.Range("Q141").Value = -0.25
$ws.Range("R141").Value = 1.925
$ws.Range("S141").Value = 1.925
$ws.Range("T141").Value = 2.75
$ws.Range("U141").Value = 2.025
$ws.Range("V141").Value = 1.825
